$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("D6").Value = "2016-32-20 20:32:44"
